$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.556.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.479.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9773"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "279.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3666"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3078"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.92"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.060"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06659"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.509"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.200"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9779"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001029"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.480.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05937"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.478"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.53"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.625.36"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.152"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.640.29"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.967"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.005"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8151"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08033"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.552"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05825"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.728"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Frax"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9772"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.736"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02046"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.47"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.529"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.22"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5199"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.800"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06466"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9954"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.51%  "
